$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 506-507 (existing rows 506..532 shift down to 508..534)
$ws.Rows("506:507").Insert()

# Row 506: new "1a amarillo" entry
$ws.Range("A506").Value = 4
$ws.Range("B506").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C506").Value = "Los Lagos"
$ws.Range("D506").Value = 44753
$ws.Range("E506").Value = 10
$ws.Range("F506").Value = "Fruta"
$ws.Range("G506").Value = 100102
$ws.Range("H506").Value = "Cítricos"
$ws.Range("I506").Value = 100102003
$ws.Range("J506").Value = "Limón"
$ws.Range("K506").Value = "Sin especificar"
$ws.Range("L506").Value = "1a amarillo"
$ws.Range("M506").Value = 500
$ws.Range("N506").Value = 9000
$ws.Range("O506").Value = 9000
$ws.Range("P506").Value = 9000
$ws.Range("Q506").Value = '$/malla 18 kilos'
$ws.Range("R506").Value = "Provincia de Melipilla"
$ws.Range("S506").Value = 500
$ws.Range("T506").Value = 18

# Row 507: new "2a amarillo" entry
$ws.Range("A507").Value = 4
$ws.Range("B507").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C507").Value = "Los Lagos"
$ws.Range("D507").Value = 44753
$ws.Range("E507").Value = 10
$ws.Range("F507").Value = "Fruta"
$ws.Range("G507").Value = 100102
$ws.Range("H507").Value = "Cítricos"
$ws.Range("I507").Value = 100102003
$ws.Range("J507").Value = "Limón"
$ws.Range("K507").Value = "Sin especificar"
$ws.Range("L507").Value = "2a amarillo"
$ws.Range("M507").Value = 250
$ws.Range("N507").Value = 8000
$ws.Range("O507").Value = 8000
$ws.Range("P507").Value = 8000
$ws.Range("Q507").Value = '$/malla 18 kilos'
$ws.Range("R507").Value = "Provincia de Melipilla"
$ws.Range("S507").Value = 444
$ws.Range("T507").Value = 18
